$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos")

# G2 -> new shared string "OSVPPRU12" (usuario column)
$ws.Range("G2").Value = "OSVPPRU12"

# H2 -> numeric value 1234 (clave column)
$ws.Range("H2").Value = 1234

# Update the active selection on the Datos sheet to F7 (matches saved file's selection)
$ws.Activate()
$ws.Range("F7").Select()
